$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Row 2 ---
$ws1.Range("A2").Value = "yes"
$ws1.Range("B2").Value = "John Fred"
$ws1.Range("C2").Value = "Firearm"
$ws1.Range("D2").Value = "dell pad"
$ws1.Range("E2").Value = "Pistol"
$ws1.Range("F2").Value = "Sativa"
$ws1.Range("G2").Value = "active"
$ws1.Range("H2").Value = "for work"
$ws1.Range("I2").Value = 45467
$ws1.Range("J2").Value = "Hadid1259"
$ws1.Range("K2").Value = "ank123559"

# --- Row 3 ---
$ws1.Range("A3").Value = "no"
$ws1.Range("B3").Value = "John Fred"
$ws1.Range("C3").Value = "Electronic Device"
$ws1.Range("D3").Value = "Iphone"
$ws1.Range("E3").Value = "Machine Gun"
$ws1.Range("F3").Value = "Sativa"
$ws1.Range("G3").Value = "active"
$ws1.Range("H3").Value = "for work"
$ws1.Range("I3").Value = 45467
$ws1.Range("J3").Value = "Ukbiased23550"
$ws1.Range("K3").Value = "mare12355"

# Match the font re-stamp seen on the date column's shared style (the style
# used by I2:I5) in the diff - the whole column style picks up the refreshed font
$ws1.Range("I2:I5").Font.Name = "Arial"

# Move the active selection to match the authored workbook state
[void]$ws1.Range("I10").Select()
